# Update cryptos list prices (column D) and volume(1h) percentages (column E)
# to reflect the latest scrape, matching the GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.985.38"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "2.594.39"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "310.55"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("D6").Value = "98.31"
$ws.Range("E6").Value = "  -2.97%  "

$ws.Range("D7").Value = "0.598"
$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").Value = "38.76"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("D11").Value = "54.33"
$ws.Range("E11").Value = "  -1.89%  "

$ws.Range("D12").Value = "0.0837"
$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("D13").Value = "8.09"
$ws.Range("E13").Value = "  -1.48%  "

$ws.Range("D14").Value = "2.996.72"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D16").Value = "2.586.66"
$ws.Range("E16").Value = "  -1.47%  "

$ws.Range("D17").Value = "0.913"
$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("D18").Value = "14.81"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("D19").Value = "46.174.97"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "12.74"
$ws.Range("E21").Value = "  -4.63%  "

$ws.Range("D22").Value = "6.68"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").Value = "293.98"
$ws.Range("E23").Value = "  +14.09%  "

$ws.Range("D24").Value = "72.73"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("D25").Value = "3.05"
$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("D27").Value = "29.56"
$ws.Range("E27").Value = "  +3.94%  "

$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("D30").Value = "10.75"
$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("D31").Value = "38.21"
$ws.Range("E31").Value = "  -5.09%  "

$ws.Range("E32").Value = "  -2.83%  "

$ws.Range("D33").Value = "6.22"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("E34").Value = "  -4.80%  "

$ws.Range("D35").Value = "155.21"
$ws.Range("E35").Value = "  +2.92%  "

$ws.Range("D36").Value = "0.0834"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").Value = "2.19"
$ws.Range("E37").Value = "  -5.74%  "

$ws.Range("D38").Value = "2.76"
$ws.Range("E38").Value = "  -6.65%  "

$ws.Range("E39").Value = "  +2.95%  "

$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").Value = "15.66"
$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").Value = "'0.0330"
$ws.Range("E42").Value = "  +1.86%  "

$ws.Range("D43").Value = "3.56"
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("D44").Value = "21.02"
$ws.Range("E44").Value = "  +9.48%  "

$ws.Range("D45").Value = "3.93"
$ws.Range("E45").Value = "  -6.08%  "

$ws.Range("D46").Value = "2.108.61"
$ws.Range("E46").Value = "  +2.96%  "

$ws.Range("D47").Value = "97.39"
$ws.Range("E47").Value = "  +6.09%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").Value = "'9.60"
$ws.Range("E49").Value = "  +3.49%  "

$ws.Range("D50").Value = "0.201"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").Value = "107.87"
$ws.Range("E51").Value = "  -1.86%  "
